# The deck's Design/Theme is switched from the custom "Integral" (Red
# Violet) palette back to the default Office Theme palette. This mirrors
# a user picking a different Theme on the Design tab: the slide master's
# ColorScheme (stored as ppt/theme/theme1.xml, the theme used by the
# slide master -> slides) is updated to the standard Office color values.
#
# VBA's ColorScheme.Colors(index).RGB uses the classic RGB() long layout
# (0x00BBGGRR), so a small helper converts plain "RRGGBB" hex bytes into
# that value.

function ToRGB($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$scheme = $master.ColorScheme

# Index order (verified empirically against the theme's <a:clrScheme>):
#  1 dk1   2 lt1   3 dk2   4 lt2
#  5 accent1  6 accent2  7 accent3  8 accent4  9 accent5  10 accent6
#  11 hlink   12 folHlink
#
# Target values = the standard Office Theme color scheme.
$scheme.Colors(1).RGB  = ToRGB 0x00 0x00 0x00   # dk1      000000
$scheme.Colors(2).RGB  = ToRGB 0xFF 0xFF 0xFF   # lt1      FFFFFF
$scheme.Colors(3).RGB  = ToRGB 0x44 0x54 0x6A   # dk2      44546A
$scheme.Colors(4).RGB  = ToRGB 0xE7 0xE6 0xE6   # lt2      E7E6E6
$scheme.Colors(5).RGB  = ToRGB 0x5B 0x9B 0xD5   # accent1  5B9BD5
$scheme.Colors(6).RGB  = ToRGB 0xED 0x7D 0x31   # accent2  ED7D31
$scheme.Colors(7).RGB  = ToRGB 0xA5 0xA5 0xA5   # accent3  A5A5A5
$scheme.Colors(8).RGB  = ToRGB 0xFF 0xC0 0x00   # accent4  FFC000
$scheme.Colors(9).RGB  = ToRGB 0x44 0x72 0xC4   # accent5  4472C4
$scheme.Colors(10).RGB = ToRGB 0x70 0xAD 0x47   # accent6  70AD47
$scheme.Colors(11).RGB = ToRGB 0x05 0x63 0xC1   # hlink    0563C1
$scheme.Colors(12).RGB = ToRGB 0x95 0x4F 0x72   # folHlink 954F72
